$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 19 - read.php: "Scroll animations?" with a reference video link in Notes
$ws.Range("A19").Value = "read.php"
$ws.Range("B19").Value = "Scroll animations?"
$ws.Range("C19").Value = "https://www.youtube.com/watch?v=0TnO1GzKWPc&ab_channel=SnippetsCode"

# Row 20 - register.php: ensure duplicate email addresses can't sign up
$ws.Range("A20").Value = "register.php"
$ws.Range("B20").Value = "Ensure you can only sign up with an email address once"
$ws.Range("C20").Value = "You'll probably need to add a clause that checks the database for that email address (if (`$row > 0) { return ""This email address is already registered.""}"

# The long wrapped note in C20 makes Excel grow row 20 to fit the text
$ws.Rows.Item(20).RowHeight = 72.5

# Mirror the author's final scroll position / selection in the sheet view
$excel.Goto($ws.Range("A16"), $false)
$ws.Range("C21").Select()
